$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168; this shifts existing rows 168..266 down to 169..267
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new record's data
$ws.Cells.Item(168, 1).Value = 8
$ws.Cells.Item(168, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(168, 3).Value = "Coquimbo"
$ws.Cells.Item(168, 4).Value = 44606
$ws.Cells.Item(168, 5).Value = 4
$ws.Cells.Item(168, 6).Value = 100112032
$ws.Cells.Item(168, 7).Value = "Zapallo italiano"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 440
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 10000
$ws.Cells.Item(168, 13).Value = 9500
$ws.Cells.Item(168, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(168, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(168, 16).Value = 158
$ws.Cells.Item(168, 17).Value = 60
$ws.Cells.Item(168, 18).Value = "Hortaliza"
